$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values in row 3 as per the commit
$ws.Range("D3").Value = 2
$ws.Range("F3").Value = -3
$ws.Range("H3").Value = 46

# Update the active selection to D3
$ws.Range("D3").Select()
